$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

function Copy-CellFormat($srcRef, $dstRef) {
    # Format-only copy so we land on the workbook's existing named cell
    # style (mtitleStyle / correctStyle / ...) instead of minting a new,
    # functionally-duplicate style entry.
    $ws.Range($srcRef).Copy()
    $ws.Range($dstRef).PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Score summary block (rows 10-12): recompute with the real marks instead of
# the placeholder "Absent" / zero values, and give the row labels the same
# bold "mtitleStyle" look as the header row above them (row 9).
# ---------------------------------------------------------------------------

Copy-CellFormat "A9" "A10"
$ws.Range("B10").Value = 25
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

Copy-CellFormat "A9" "A11"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

Copy-CellFormat "A9" "A12"
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "100/112"

# ---------------------------------------------------------------------------
# The per-question table used to have three Student-Ans/Correct-Ans column
# pairs (A-B, D-E, G-H). Only the first pair is actually needed, so drop the
# third pair entirely and trim the second pair back to just the header rows
# that stay populated (16-18).
# ---------------------------------------------------------------------------

$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# Fill in the remaining second-pair "Student Ans" cells (D16:D18) to mirror
# their "Correct Ans" counterparts, using the same "correctStyle" look that
# the first-pair Student Ans column below will get.
Copy-CellFormat "B10" "D16"
$ws.Range("D16").Value = "Option A"

Copy-CellFormat "B10" "D17"
$ws.Range("D17").Value = "Option C"

Copy-CellFormat "B10" "D18"
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------------
# Column A ("Student Ans") for the first pair was blank for every question;
# populate it with the student's actual answers now that they're known.
# Rows 28 and 29 remain unattempted, so they stay blank.
# ---------------------------------------------------------------------------

$studentAnswers = [ordered]@{
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    25 = "Option A"
    26 = "Option C"
    27 = "Option A"
    30 = "Option B"
    31 = "Option D"
    32 = "Option C"
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

foreach ($row in $studentAnswers.Keys) {
    Copy-CellFormat "B10" "A$row"
    $ws.Range("A$row").Value = $studentAnswers[$row]
}
